$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (bottom-up to keep row numbers stable):
# Row 27: Mystery Booster Playtest Cards -> Mystery Booster Playtest Cards 2021
# Row 26: Promo Pack: Streets of New Capenna -> Streets of New Capenna Promos
# Row 7:  Mystery Booster Playtest Cards 2021 -> Mystery Booster Playtest Cards
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(7).Delete()

$ws.Range("B14").Select()
